$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C8").Value = 20000
$ws.Range("C9").Value = 30000
$ws.Range("C10").Value = 35000
$ws.Range("C11").Value = 25000
$ws.Range("C12").Value = 20000
$ws.Range("C13").Value = 15000

$ws.Range("C15").Value = 18000
$ws.Range("C16").Value = 25000
$ws.Range("C17").Value = 20000

$ws.Range("C20").Value = 15000
$ws.Range("C21").Value = 25000
$ws.Range("C22").Value = 23000
$ws.Range("C23").Value = 20000

$ws.Range("F21").Select()
